$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Week 3 time log entries (rows 24-27)
$ws.Range("B24").Value = 10
$ws.Range("C24").Value = 10

$ws.Range("B25").Value = 20
$ws.Range("C25").Value = 15

$ws.Range("B26").Value = 120
$ws.Range("C26").Value = 50

$ws.Range("B27").Value = 20
$ws.Range("C27").Value = 10

# Update the visible scroll position / selection to match the author's
# working state when they finished this edit.
$ws.Application.Goto($ws.Range("A15"), $true)
$ws.Range("D27").Select()
